# Insert a new data row at row 121 (pushing existing rows 121-244 down to 122-245)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(121).Insert()

$ws.Range("A121").Value = 10
$ws.Range("B121").Value = "Vega Modelo de Temuco"
$ws.Range("C121").Value = "La Araucanía"
$ws.Range("D121").Value = 44586
$ws.Range("E121").Value = 9
$ws.Range("F121").Value = 100112017
$ws.Range("G121").Value = "Apio"
$ws.Range("H121").Value = "Americana (o)"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 65
$ws.Range("K121").Value = 10000
$ws.Range("L121").Value = 10000
$ws.Range("M121").Value = 10000
$ws.Range("N121").Value = "$/docena de matas"
$ws.Range("O121").Value = "Provincia del Elquí"
$ws.Range("P121").Value = 1667
$ws.Range("Q121").Value = 6
$ws.Range("R121").Value = "Hortaliza"
